# Weekly fruit/vegetable price update: insert 3 new Chirimoya price rows
# (Macroferia Regional de Talca) at the top of the data block, pushing the
# existing rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 83 (shifts old rows 83:104 -> 86:107,
# and grows the sheet dimension to A1:T107 automatically).
$ws.Range("A83:T85").Insert()

# Values that are constant across every data row in this sheet.
$marketId = 5
$market = "Macroferia Regional de Talca"
$region = "Maule"
$codreg = 7
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 10 kilos"
$origen = "Provincia de Limarí"
$kgUnidad = 10

function Set-Row($r, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($r, 1).Value = $marketId
    $ws.Cells.Item($r, 2).Value = $market
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 83 44841 "Especial" 60 22000 22000 22000 2200
Set-Row 84 44841 "Extra (doble especial)" 50 25000 25000 25000 2500
Set-Row 85 44841 "Primera" 40 20000 20000 20000 2000
